$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100: mark as Committed, add Code Review + Committed dates, add revision number
$ws.Range("D100").Value = "Committed"
$ws.Range("E100").Copy($ws.Range("F100"))
$ws.Range("E100").Copy($ws.Range("G100"))
$ws.Range("D100").Copy($ws.Range("H100"))
$ws.Range("H100").Value = "5edb7be"

# Row 101: new entry
$ws.Range("A99").Copy($ws.Range("A101"))
$ws.Range("A101").Value = "Created NC_EDGAR_sector_mapping.csv and Master_EDGAR_sector_mapping.csv for use in module C and checking presence of all EDGAR sectors in final emissions database."
$ws.Range("B99").Copy($ws.Range("B101"))
$ws.Range("B101").Value = 97
$ws.Range("C99").Copy($ws.Range("C101"))
$ws.Range("D99").Copy($ws.Range("D101"))
$ws.Range("D101").Value = "Review"
$ws.Range("E99").Copy($ws.Range("E101"))
$ws.Range("E101").Value = 42380
$ws.Rows("101").RowHeight = 62.25

# Row 102: new entry
$ws.Range("A99").Copy($ws.Range("A102"))
$ws.Range("A102").Value = "Added EDGARcheck function to analysis_functions.R"
$ws.Range("B99").Copy($ws.Range("B102"))
$ws.Range("B102").Value = 98
$ws.Range("C99").Copy($ws.Range("C102"))
$ws.Range("D99").Copy($ws.Range("D102"))
$ws.Range("D102").Value = "Review"
$ws.Range("E99").Copy($ws.Range("E102"))
$ws.Range("E102").Value = 42380
$ws.Rows("102").RowHeight = 32.25

[void]$ws.Range("E107").Select()
